# Day 19 and 20 Solutions and runtimes
# Fill in the Part 1 / Part 2 run times for Day 20 (row 24) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B24").Value = 0.23102749988902299
$ws.Range("C24").Value = 0.00084280001465231104

$excel.CalculateFullRebuild()
